# Restructure the worksheet: add a numeric column-index row, split the header
# row out from row 1, factor the constant "Polycarbonate Plastic" material value
# into its own row, and drop the now-redundant per-row material_surface values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing header row (and all data beneath it) down by two rows so
# we can introduce the new index row and the material-note row above it.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(3).Insert()

# New row 1: plain numeric column indexes (0-11), keeps the bold/border style
# that used to belong to the header row.
$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 7
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 9
$ws.Cells.Item(1, 11).Value = 10
$ws.Cells.Item(1, 12).Value = 11

# New row 2: the original column headers, now unstyled.
$ws.Cells.Item(2, 1).Value = 'Lg.'
$ws.Cells.Item(2, 2).Value = 'Threading'
$ws.Cells.Item(2, 3).Value = 'HeadDia.'
$ws.Cells.Item(2, 4).Value = 'Head Ht.'
$ws.Cells.Item(2, 5).Value = 'DriveSize'
$ws.Cells.Item(2, 6).Value = 'TemperatureRange, °F'
$ws.Cells.Item(2, 7).Value = 'Color'
$ws.Cells.Item(2, 8).Value = 'Pkg.Qty.'
$ws.Cells.Item(2, 10).Value = 'Pkg.'
$ws.Range("A2:L2").ClearFormats()

# New row 3: a single note cell carrying the material literal, unstyled.
$ws.Cells.Item(3, 1).Value = 'Polycarbonate Plastic'
$ws.Range("A3:L3").ClearFormats()

# Every data row (now at rows 4-34) keeps its own values untouched, except
# column L (material_surface), whose constant value moved to A3 above.
for ($r = 4; $r -le 34; $r++) {
    $ws.Cells.Item($r, 12).ClearContents()
}
